$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reconciliation function adjustment: column O truncated to date-only,
# BP/BQ status columns corrected (NaoIniciado=716 / Concluido=784),
# CF status-bar counter shifted by +264.
$ws.Range("O21").Value = 45337
$ws.Range("BP21").Value = "Concluido"
$ws.Range("BQ21").Value = "NaoIniciado"
$ws.Range("CF21").Value = 281
$ws.Range("O22").Value = 45337
$ws.Range("BP22").Value = "NaoIniciado"
$ws.Range("BQ22").Value = "Concluido"
$ws.Range("CF22").Value = 282
$ws.Range("O23").Value = 45337
$ws.Range("BP23").Value = "Concluido"
$ws.Range("BQ23").Value = "NaoIniciado"
$ws.Range("CF23").Value = 281
$ws.Range("O24").Value = 45337
$ws.Range("BP24").Value = "NaoIniciado"
$ws.Range("BQ24").Value = "Concluido"
$ws.Range("CF24").Value = 282
$ws.Range("O25").Value = 45337
$ws.Range("BP25").Value = "NaoIniciado"
$ws.Range("BQ25").Value = "Concluido"
$ws.Range("CF25").Value = 282
$ws.Range("O26").Value = 45337
$ws.Range("BP26").Value = "Concluido"
$ws.Range("BQ26").Value = "NaoIniciado"
$ws.Range("CF26").Value = 283
$ws.Range("O27").Value = 45337
$ws.Range("BP27").Value = "Concluido"
$ws.Range("BQ27").Value = "NaoIniciado"
$ws.Range("CF27").Value = 282
$ws.Range("O28").Value = 45337
$ws.Range("BP28").Value = "Concluido"
$ws.Range("BQ28").Value = "NaoIniciado"
$ws.Range("CF28").Value = 283
$ws.Range("O29").Value = 45337
$ws.Range("BP29").Value = "Concluido"
$ws.Range("BQ29").Value = "NaoIniciado"
$ws.Range("CF29").Value = 283
$ws.Range("O30").Value = 45337
$ws.Range("BP30").Value = "NaoIniciado"
$ws.Range("BQ30").Value = "Concluido"
$ws.Range("CF30").Value = 283
$ws.Range("O31").Value = 45337
$ws.Range("BP31").Value = "NaoIniciado"
$ws.Range("BQ31").Value = "Concluido"
$ws.Range("CF31").Value = 283
$ws.Range("O32").Value = 45337
$ws.Range("BP32").Value = "Concluido"
$ws.Range("BQ32").Value = "NaoIniciado"
$ws.Range("CF32").Value = 282
$ws.Range("O33").Value = 45337
$ws.Range("BP33").Value = "Concluido"
$ws.Range("BQ33").Value = "NaoIniciado"
$ws.Range("CF33").Value = 282
$ws.Range("O34").Value = 45337
$ws.Range("BP34").Value = "Concluido"
$ws.Range("BQ34").Value = "NaoIniciado"
$ws.Range("CF34").Value = 282
$ws.Range("O35").Value = 45337
$ws.Range("BP35").Value = "NaoIniciado"
$ws.Range("BQ35").Value = "Concluido"
$ws.Range("CF35").Value = 282
$ws.Range("O36").Value = 45337
$ws.Range("BP36").Value = "NaoIniciado"
$ws.Range("BQ36").Value = "Concluido"
$ws.Range("CF36").Value = 282
$ws.Range("O37").Value = 45337
$ws.Range("BP37").Value = "Concluido"
$ws.Range("BQ37").Value = "NaoIniciado"
$ws.Range("CF37").Value = 283
$ws.Range("O38").Value = 45337
$ws.Range("BP38").Value = "NaoIniciado"
$ws.Range("BQ38").Value = "Concluido"
$ws.Range("CF38").Value = 283
$ws.Range("O39").Value = 45337
$ws.Range("BP39").Value = "Concluido"
$ws.Range("BQ39").Value = "NaoIniciado"
$ws.Range("CF39").Value = 283
$ws.Range("O40").Value = 45337
$ws.Range("BP40").Value = "Concluido"
$ws.Range("BQ40").Value = "NaoIniciado"
$ws.Range("CF40").Value = 283
$ws.Range("O41").Value = 45337
$ws.Range("BP41").Value = "NaoIniciado"
$ws.Range("BQ41").Value = "Concluido"
$ws.Range("CF41").Value = 282
$ws.Range("O42").Value = 45337
$ws.Range("BP42").Value = "NaoIniciado"
$ws.Range("BQ42").Value = "Concluido"
$ws.Range("CF42").Value = 285
$ws.Range("O43").Value = 45337
$ws.Range("BP43").Value = "NaoIniciado"
$ws.Range("BQ43").Value = "Concluido"
$ws.Range("CF43").Value = 283
$ws.Range("O44").Value = 45337
$ws.Range("BP44").Value = "Concluido"
$ws.Range("BQ44").Value = "NaoIniciado"
$ws.Range("CF44").Value = 283
$ws.Range("O45").Value = 45337
$ws.Range("BP45").Value = "Concluido"
$ws.Range("BQ45").Value = "NaoIniciado"
$ws.Range("CF45").Value = 285
$ws.Range("O46").Value = 45337
$ws.Range("BP46").Value = "Concluido"
$ws.Range("BQ46").Value = "NaoIniciado"
$ws.Range("CF46").Value = 285
$ws.Range("O47").Value = 45337
$ws.Range("BP47").Value = "NaoIniciado"
$ws.Range("BQ47").Value = "Concluido"
$ws.Range("CF47").Value = 283
$ws.Range("O48").Value = 45337
$ws.Range("BP48").Value = "NaoIniciado"
$ws.Range("BQ48").Value = "Concluido"
$ws.Range("CF48").Value = 283
$ws.Range("O49").Value = 45337
$ws.Range("BP49").Value = "NaoIniciado"
$ws.Range("BQ49").Value = "Concluido"
$ws.Range("CF49").Value = 283
$ws.Range("O50").Value = 45337
$ws.Range("BP50").Value = "NaoIniciado"
$ws.Range("BQ50").Value = "Concluido"
$ws.Range("CF50").Value = 283
$ws.Range("O51").Value = 45337
$ws.Range("BP51").Value = "Concluido"
$ws.Range("BQ51").Value = "NaoIniciado"
$ws.Range("CF51").Value = 285
$ws.Range("O52").Value = 45337
$ws.Range("BP52").Value = "NaoIniciado"
$ws.Range("BQ52").Value = "Concluido"
$ws.Range("CF52").Value = 283
$ws.Range("O53").Value = 45337
$ws.Range("BP53").Value = "NaoIniciado"
$ws.Range("BQ53").Value = "Concluido"
$ws.Range("CF53").Value = 285
$ws.Range("O54").Value = 45337
$ws.Range("BP54").Value = "NaoIniciado"
$ws.Range("BQ54").Value = "Concluido"
$ws.Range("CF54").Value = 284
$ws.Range("O55").Value = 45337
$ws.Range("BP55").Value = "Concluido"
$ws.Range("BQ55").Value = "NaoIniciado"
$ws.Range("CF55").Value = 286
$ws.Range("O56").Value = 45337
$ws.Range("BP56").Value = "Concluido"
$ws.Range("BQ56").Value = "NaoIniciado"
$ws.Range("CF56").Value = 283
$ws.Range("O57").Value = 45337
$ws.Range("BP57").Value = "Concluido"
$ws.Range("BQ57").Value = "NaoIniciado"
$ws.Range("CF57").Value = 286
$ws.Range("O58").Value = 45337
$ws.Range("BP58").Value = "NaoIniciado"
$ws.Range("BQ58").Value = "Concluido"
$ws.Range("CF58").Value = 287
$ws.Range("O59").Value = 45337
$ws.Range("BP59").Value = "Concluido"
$ws.Range("BQ59").Value = "NaoIniciado"
$ws.Range("CF59").Value = 283
$ws.Range("O60").Value = 45337
$ws.Range("BP60").Value = "Concluido"
$ws.Range("BQ60").Value = "NaoIniciado"
$ws.Range("CF60").Value = 285
$ws.Range("O61").Value = 45337
$ws.Range("BP61").Value = "Concluido"
$ws.Range("BQ61").Value = "NaoIniciado"
$ws.Range("CF61").Value = 287
$ws.Range("O62").Value = 45337
$ws.Range("BP62").Value = "Concluido"
$ws.Range("BQ62").Value = "NaoIniciado"
$ws.Range("CF62").Value = 288
$ws.Range("O63").Value = 45337
$ws.Range("BP63").Value = "Concluido"
$ws.Range("BQ63").Value = "NaoIniciado"
$ws.Range("CF63").Value = 287
$ws.Range("O64").Value = 45337
$ws.Range("BP64").Value = "Concluido"
$ws.Range("BQ64").Value = "NaoIniciado"
$ws.Range("CF64").Value = 288
$ws.Range("O65").Value = 45337
$ws.Range("BP65").Value = "Concluido"
$ws.Range("BQ65").Value = "NaoIniciado"
$ws.Range("CF65").Value = 287
